$wb = $excel.ActiveWorkbook

# --- 1. Un-select current active tab (AFR 3) so the new sheet becomes active ---
$afr3 = $wb.Worksheets.Item("AFR 3 (RX8 Evolve 2006 UK)")

# --- 2. Add the new worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Volumetry Efficienccy "

# --- 3. Populate the Volumetric Efficiency table (A1:R20) ---
$veData = @(
  @(1.0792999999999999,1.0792999999999999,1.0556000000000001,1.0704,1.0382,1.0286,1.0496000000000001,1.0505,0.99950000000000006,0.97360000000000002,0.96150000000000002,1.0179,1.1039000000000001,0.99329999999999996,0.99870000000000003,0.99580000000000002,1.0017,1.0154000000000001),
  @(1.0792999999999999,1.0792999999999999,1.0596000000000001,1.0205,1.0149999999999999,1.0293000000000001,1.0368999999999999,1.022,0.96430000000000005,0.97130000000000005,0.96409999999999996,1.0316000000000001,1.1039000000000001,0.99329999999999996,0.99870000000000003,0.99580000000000002,1.0017,1.0154000000000001),
  @(1.0616000000000001,1.0616000000000001,1.0556000000000001,1.0125999999999999,1.0095000000000001,1.0165,1.0367999999999999,1.0135000000000001,0.96679999999999999,0.96909999999999996,0.9667,1.0452999999999999,1.0504,0.99329999999999996,0.99870000000000003,0.99580000000000002,1.0017,1.0154000000000001),
  @(1.0609,1.0609,1.0295000000000001,1.0024,1.0009999999999999,1.0135000000000001,1.0116000000000001,0.99650000000000005,0.97640000000000005,0.9556,0.96230000000000004,1.0219,0.997,1.0277000000000001,1.0075000000000001,0.99580000000000002,0.99080000000000001,0.99580000000000002),
  @(1.0555000000000001,1.0555000000000001,1.0176000000000001,0.99539999999999995,0.97760000000000002,1.004,0.98050000000000004,0.98429999999999995,0.97250000000000003,0.94210000000000005,0.9698,0.9839,1.0055000000000001,1.0474000000000001,1.0105,0.99080000000000001,0.99580000000000002,0.98099999999999998),
  @(1.0532999999999999,1.0532999999999999,0.99729999999999996,0.97460000000000002,0.96989999999999998,0.98450000000000004,0.96540000000000004,0.96260000000000001,0.96860000000000002,0.94899999999999995,0.94650000000000001,1.0122,1.0246,1.0498000000000001,0.99580000000000002,0.9879,0.99580000000000002,1.0203),
  @(1.0528,1.0528,0.99270000000000003,0.96750000000000003,0.9677,0.9708,0.9617,0.94099999999999995,0.95530000000000004,0.94710000000000005,0.93289999999999995,1.0148999999999999,1.0106999999999999,1.0251999999999999,0.98,1.0056,0.9889,1.0056),
  @(1.0291999999999999,1.0291999999999999,0.97850000000000004,0.97219999999999995,0.97089999999999999,0.96179999999999999,0.95599999999999996,0.95209999999999995,0.94369999999999998,0.9325,0.97709999999999997,1.0119,1.0007999999999999,0.99329999999999996,1.0056,0.99580000000000002,0.9859,0.99580000000000002),
  @(1.0056,1.0056,0.96430000000000005,0.95779999999999998,0.97970000000000002,0.96179999999999999,0.95030000000000003,0.94869999999999999,0.93669999999999998,0.97689999999999999,0.98099999999999998,1.0002,0.99080000000000001,1.0326,0.99580000000000002,0.99970000000000003,0.98399999999999999,0.97609999999999997),
  @(0.99580000000000002,0.99580000000000002,0.96909999999999996,0.94330000000000003,0.96309999999999996,0.94269999999999998,0.95309999999999995,0.94450000000000001,0.95109999999999995,0.9597,0.98960000000000004,0.98839999999999995,0.99080000000000001,1.0179,0.9879,0.97899999999999998,0.97219999999999995,0.96630000000000005),
  @(0.9859,0.9859,0.95150000000000001,0.94269999999999998,0.9466,0.92359999999999998,0.93240000000000001,0.93610000000000004,0.95209999999999995,0.97030000000000005,0.99819999999999998,0.96630000000000005,0.99329999999999996,0.99329999999999996,0.98099999999999998,0.99080000000000001,0.97609999999999997,0.96630000000000005),
  @(0.9859,0.9859,0.95150000000000001,0.93189999999999995,0.9073,0.92200000000000004,0.9133,0.92579999999999996,0.93510000000000004,0.98099999999999998,0.98350000000000004,0.9859,0.98099999999999998,0.97860000000000003,0.96630000000000005,0.98399999999999999,0.96630000000000005,0.9466),
  @(0.9859,0.9859,0.95150000000000001,0.93189999999999995,0.9073,0.90239999999999998,0.88619999999999999,0.91020000000000001,0.92920000000000003,0.99819999999999998,0.99580000000000002,0.97860000000000003,0.96870000000000001,0.96630000000000005,0.95640000000000003,0.97219999999999995,0.95150000000000001,0.95840000000000003),
  @(0.9859,0.9859,0.95150000000000001,0.93189999999999995,0.9073,0.90239999999999998,0.90590000000000004,0.93430000000000002,0.97609999999999997,0.98099999999999998,0.97119999999999995,0.96630000000000005,0.94910000000000005,0.96870000000000001,0.95640000000000003,0.99580000000000002,0.97319999999999995,0.97609999999999997),
  @(0.9859,0.9859,0.95150000000000001,0.93189999999999995,0.9073,0.90239999999999998,0.90590000000000004,0.93430000000000002,0.97609999999999997,0.96060000000000001,0.97119999999999995,0.9627,0.94189999999999996,0.96509999999999996,0.97299999999999998,0.99580000000000002,0.97319999999999995,0.89390000000000003),
  @(0.9859,0.9859,0.95150000000000001,0.93189999999999995,0.9073,0.90239999999999998,0.90590000000000004,0.93430000000000002,0.97609999999999997,0.95960000000000001,0.94740000000000002,0.96140000000000003,0.93920000000000003,0.95450000000000002,0.97899999999999998,0.96479999999999999,0.93510000000000004,0.89190000000000003),
  @(0.9859,0.9859,0.95150000000000001,0.93189999999999995,0.9073,0.90239999999999998,0.90590000000000004,0.93430000000000002,0.97609999999999997,0.95850000000000002,0.92200000000000004,0.94699999999999995,0.92479999999999996,0.94320000000000004,0.97899999999999998,0.93189999999999995,0.89449999999999996,0.88970000000000005),
  @(0.9859,0.9859,0.95150000000000001,0.93189999999999995,0.9073,0.90239999999999998,0.90590000000000004,0.93430000000000002,0.97609999999999997,0.95589999999999997,0.91549999999999998,0.91259999999999997,0.89049999999999996,0.91610000000000003,0.92989999999999995,0.9254,0.88800000000000001,0.88460000000000005),
  @(0.9859,0.9859,0.95150000000000001,0.93189999999999995,0.9073,0.90239999999999998,0.90590000000000004,0.93430000000000002,0.97609999999999997,0.95069999999999999,0.90239999999999998,0.90080000000000005,0.87870000000000004,0.90429999999999999,0.91810000000000003,0.91220000000000001,0.87490000000000001,0.87419999999999998),
  @(0.9859,0.9859,0.95150000000000001,0.93189999999999995,0.9073,0.90239999999999998,0.90590000000000004,0.93430000000000002,0.97609999999999997,0.95069999999999999,0.90239999999999998,0.89290000000000003,0.87080000000000002,0.89649999999999996,0.91020000000000001,0.91220000000000001,0.87490000000000001,0.87419999999999998)
)

for ($r = 0; $r -lt $veData.Count; $r++) {
  $rowVals = $veData[$r]
  for ($c = 0; $c -lt $rowVals.Count; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
  }
}

# --- 4. Set selection / view on new sheet to match target (D25 selected, no zoom overrides) ---
$ws.Range("D25").Select()

# --- 5. Activate the new sheet so it becomes the selected tab ---
$ws.Activate()
